# porcentajes-participacion.xlsx - add "Hoja2" (marcas, etiquetas y productos)
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Update Hoja1's view: drop tabSelected, zoom to 70%, select the whole used range ---
$ws1.Range("A1:E21").Select()

# --- Add the new worksheet right after Hoja1 ---
$ws = $wb.Worksheets.Add($null, $ws1)
$ws.Name = "Hoja2"

# Header row
$ws.Range("A1").Value = "Grupo"
$ws.Range("B1").Value = "Responsable"
$ws.Range("C1").Value = "Alumno"
$ws.Range("D1").Value = "Porcentaje"
$ws.Range("E1").Value = "Justificación"

# ISOs group
$ws.Range("A2").Value = "ISOs"
$ws.Range("B2").Value = "Rivera La Rosa, Jaime Jacob"
$ws.Range("C2").Value = "Calcina Pacori, Julio Rolando"
$ws.Range("D2").Value = 0.3
$ws.Range("D2").NumberFormat = "0%"

$ws.Range("A3").Value = "ISOs"
$ws.Range("B3").Value = "Rivera La Rosa, Jaime Jacob"
$ws.Range("C3").Value = "Marin Nuñez, José"
$ws.Range("D3").Value = 0.3
$ws.Range("D3").NumberFormat = "0%"

$ws.Range("A4").Value = "ISOs"
$ws.Range("B4").Value = "Rivera La Rosa, Jaime Jacob"
$ws.Range("C4").Value = "León Crispín, Mao Jovaldo"
$ws.Range("D4").Value = 0.3
$ws.Range("D4").NumberFormat = "0%"

$ws.Range("A5").Value = "ISOs"
$ws.Range("B5").Value = "Rivera La Rosa, Jaime Jacob"
$ws.Range("C5").Value = "Silva Dolores, Cristian"
$ws.Range("D5").Value = 0.3
$ws.Range("D5").NumberFormat = "0%"

$ws.Range("A6").Value = "ISOs"
$ws.Range("B6").Value = "Rivera La Rosa, Jaime Jacob"
$ws.Range("C6").Value = "Yánac Jiménes Daniel Nehemías"
$ws.Range("D6").Value = 0.3
$ws.Range("D6").NumberFormat = "0%"

$ws.Range("A7").Value = "ISOs"
$ws.Range("B7").Value = "Rivera La Rosa, Jaime Jacob"
$ws.Range("C7").Value = "Montoro Correa, Jesús Alberto"
$ws.Range("D7").Value = 0.3
$ws.Range("D7").NumberFormat = "0%"

# Desarrollo group
$ws.Range("A8").Value = "Desarrollo"
$ws.Range("B8").Value = "Luis Mendoza, Samuel"
$ws.Range("C8").Value = "Sandon Mateo Amilcar"
$ws.Range("D8").Value = 0.3
$ws.Range("D8").NumberFormat = "0%"

$ws.Range("A9").Value = "Desarrollo"
$ws.Range("B9").Value = "Luis Mendoza, Samuel"
$ws.Range("C9").Value = "Solis Ocaña, Luis"
$ws.Range("D9").Value = 0.3
$ws.Range("D9").NumberFormat = "0%"

$ws.Range("A10").Value = "Desarrollo"
$ws.Range("B10").Value = "Luis Mendoza, Samuel"
$ws.Range("C10").Value = "Valenzuela Paucar, Junior"
$ws.Range("D10").Value = 0.3
$ws.Range("D10").NumberFormat = "0%"

$ws.Range("A11").Value = "Desarrollo"
$ws.Range("B11").Value = "Luis Mendoza, Samuel"
$ws.Range("C11").Value = "Príncipe Henostroza, Jhordy Anderson"
$ws.Range("D11").Value = 0.5
$ws.Range("D11").NumberFormat = "0%"

$ws.Range("A12").Value = "Desarrollo"
$ws.Range("B12").Value = "Luis Mendoza, Samuel"
$ws.Range("C12").Value = "Meléndez Panana, César Yair"
$ws.Range("D12").Value = 1
$ws.Range("D12").NumberFormat = "0%"

$ws.Range("A13").Value = "Desarrollo"
$ws.Range("B13").Value = "Luis Mendoza, Samuel"
$ws.Range("C13").Value = "Cabrel Espinoza, Luis Alejandro"
$ws.Range("D13").Value = 0.3
$ws.Range("D13").NumberFormat = "0%"

$ws.Range("A14").Value = "Desarrollo"
$ws.Range("B14").Value = "Luis Mendoza, Samuel"
$ws.Range("C14").Value = "Guerrero Sanchez, Weiner Brayan"
$ws.Range("D14").Value = 0.5
$ws.Range("D14").NumberFormat = "0%"

$ws.Range("A15").Value = "Desarrollo"
$ws.Range("B15").Value = "Luis Mendoza, Samuel"
$ws.Range("C15").Value = "Porlles Pardo, Julio Hagi"
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = "0%"

$ws.Range("A16").Value = "Desarrollo"
$ws.Range("B16").Value = "Luis Mendoza, Samuel"
$ws.Range("C16").Value = "López Castro, Angello Jones"
$ws.Range("D16").Value = 1
$ws.Range("D16").NumberFormat = "0%"

# Footer block (rows 20-24), centered and merged across A:D
$ws.Range("A20").Value = "OCTAVO CICLO "
$ws.Range("A21").Value = "Ingeniería Informática"
$ws.Range("A22").Value = "Universidad Nacional José Faustino Sánchez Carrión"
$ws.Range("A23").Value = "Huacho"
$ws.Range("A24").Value = "15 de julio del 2019"

$ws.Range("A20:D20").HorizontalAlignment = -4108
$ws.Range("A21:D21").HorizontalAlignment = -4108
$ws.Range("A22:D22").HorizontalAlignment = -4108
$ws.Range("A23:D23").HorizontalAlignment = -4108
$ws.Range("A24:D24").HorizontalAlignment = -4108

$ws.Range("A20:D20").Merge()
$ws.Range("A21:D21").Merge()
$ws.Range("A22:D22").Merge()
$ws.Range("A23:D23").Merge()
$ws.Range("A24:D24").Merge()

# Column widths (character units) to reproduce the target OOXML widths as closely as possible
$ws.Columns.Item(1).ColumnWidth = 9.833333333333332
$ws.Columns.Item(2).ColumnWidth = 27.833333333333332
$ws.Columns.Item(3).ColumnWidth = 37.33333333333333
$ws.Columns.Item(5).ColumnWidth = 237.66666666666666

# Make Hoja2 the active sheet/tab, zoomed to 70%, cursor on E4
$ws.Activate()
$ws.Range("E4").Select()
$excel.ActiveWindow.Zoom = 70

$ws1.Activate()
$excel.ActiveWindow.Zoom = 70
$ws.Activate()
